$d = $word.ActiveDocument
$tbl = $d.Tables(1)

# Row 1
$tbl.Cell(1, 1).Range.Text = "53÷8=6, 5"
$tbl.Cell(1, 2).Range.Text = "66÷9=7, 3"
$tbl.Cell(1, 3).Range.Text = "97÷7=13, 6"
$tbl.Cell(1, 4).Range.Text = "70÷3=23, 1"
$tbl.Cell(1, 5).Range.Text = "44÷3=14, 2"

# Row 5
$tbl.Cell(5, 1).Range.Text = "32÷3=10, 2"
$tbl.Cell(5, 2).Range.Text = "13÷7=1, 6"
$tbl.Cell(5, 3).Range.Text = "99÷4=24, 3"
$tbl.Cell(5, 4).Range.Text = "89÷7=12, 5"
$tbl.Cell(5, 5).Range.Text = "62÷3=20, 2"

# Row 9
$tbl.Cell(9, 1).Range.Text = "89÷4=22, 1"
$tbl.Cell(9, 2).Range.Text = "90÷4=22, 2"
$tbl.Cell(9, 3).Range.Text = "77÷6=12, 5"
$tbl.Cell(9, 4).Range.Text = "59÷4=14, 3"
$tbl.Cell(9, 5).Range.Text = "10÷5=2, 0"

# Row 13
$tbl.Cell(13, 1).Range.Text = "11÷4=2, 3"
$tbl.Cell(13, 2).Range.Text = "63÷9=7, 0"
$tbl.Cell(13, 3).Range.Text = "86÷9=9, 5"
$tbl.Cell(13, 4).Range.Text = "42÷4=10, 2"
$tbl.Cell(13, 5).Range.Text = "10÷8=1, 2"

# Row 17
$tbl.Cell(17, 1).Range.Text = "38÷5=7, 3"
$tbl.Cell(17, 2).Range.Text = "67÷6=11, 1"
$tbl.Cell(17, 3).Range.Text = "64÷2=32, 0"
$tbl.Cell(17, 4).Range.Text = "49÷7=7, 0"
$tbl.Cell(17, 5).Range.Text = "45÷3=15, 0"
